$d = $word.ActiveDocument

$d.Content.Find.Execute("44-16=", $true, $false, $false, $false, $false, $true, 1, $false, "51+3=", 2) | Out-Null
$d.Content.Find.Execute("11+72=", $true, $false, $false, $false, $false, $true, 1, $false, "15-13=", 2) | Out-Null
$d.Content.Find.Execute("68+25=", $true, $false, $false, $false, $false, $true, 1, $false, "50-22=", 2) | Out-Null
$d.Content.Find.Execute("39-38=", $true, $false, $false, $false, $false, $true, 1, $false, "73-58=", 2) | Out-Null
$d.Content.Find.Execute("98-47=", $true, $false, $false, $false, $false, $true, 1, $false, "60+11=", 2) | Out-Null
$d.Content.Find.Execute("34+9=", $true, $false, $false, $false, $false, $true, 1, $false, "10+62=", 2) | Out-Null
$d.Content.Find.Execute("95-83=", $true, $false, $false, $false, $false, $true, 1, $false, "58-23=", 2) | Out-Null
$d.Content.Find.Execute("70-33=", $true, $false, $false, $false, $false, $true, 1, $false, "40+20=", 2) | Out-Null
$d.Content.Find.Execute("94-22=", $true, $false, $false, $false, $false, $true, 1, $false, "13+67=", 2) | Out-Null
$d.Content.Find.Execute("43+18=", $true, $false, $false, $false, $false, $true, 1, $false, "74-15=", 2) | Out-Null
$d.Content.Find.Execute("92-1=", $true, $false, $false, $false, $false, $true, 1, $false, "42+38=", 2) | Out-Null
$d.Content.Find.Execute("68-59=", $true, $false, $false, $false, $false, $true, 1, $false, "72-15=", 2) | Out-Null
$d.Content.Find.Execute("97-68=", $true, $false, $false, $false, $false, $true, 1, $false, "67-10=", 2) | Out-Null
$d.Content.Find.Execute("76-10=", $true, $false, $false, $false, $false, $true, 1, $false, "10+11=", 2) | Out-Null
$d.Content.Find.Execute("56+2=", $true, $false, $false, $false, $false, $true, 1, $false, "87-65=", 2) | Out-Null
$d.Content.Find.Execute("61+4=", $true, $false, $false, $false, $false, $true, 1, $false, "26-22=", 2) | Out-Null
$d.Content.Find.Execute("39+60=", $true, $false, $false, $false, $false, $true, 1, $false, "5+37=", 2) | Out-Null
$d.Content.Find.Execute("11+88=", $true, $false, $false, $false, $false, $true, 1, $false, "19-11=", 2) | Out-Null
$d.Content.Find.Execute("23-15=", $true, $false, $false, $false, $false, $true, 1, $false, "67+18=", 2) | Out-Null
$d.Content.Find.Execute("37+45=", $true, $false, $false, $false, $false, $true, 1, $false, "10+4=", 2) | Out-Null
$d.Content.Find.Execute("87+9=", $true, $false, $false, $false, $false, $true, 1, $false, "89-15=", 2) | Out-Null
$d.Content.Find.Execute("97-3=", $true, $false, $false, $false, $false, $true, 1, $false, "29+55=", 2) | Out-Null
$d.Content.Find.Execute("61-51=", $true, $false, $false, $false, $false, $true, 1, $false, "36+11=", 2) | Out-Null
$d.Content.Find.Execute("21+44=", $true, $false, $false, $false, $false, $true, 1, $false, "55+38=", 2) | Out-Null
$d.Content.Find.Execute("18+53=", $true, $false, $false, $false, $false, $true, 1, $false, "94-67=", 2) | Out-Null
$d.Content.Find.Execute("72+3=", $true, $false, $false, $false, $false, $true, 1, $false, "1+56=", 2) | Out-Null
$d.Content.Find.Execute("87-87=", $true, $false, $false, $false, $false, $true, 1, $false, "70-5=", 2) | Out-Null
$d.Content.Find.Execute("91-26=", $true, $false, $false, $false, $false, $true, 1, $false, "26+60=", 2) | Out-Null
$d.Content.Find.Execute("82-0=", $true, $false, $false, $false, $false, $true, 1, $false, "70+8=", 2) | Out-Null
$d.Content.Find.Execute("98-74=", $true, $false, $false, $false, $false, $true, 1, $false, "89+0=", 2) | Out-Null
$d.Content.Find.Execute("41+54=", $true, $false, $false, $false, $false, $true, 1, $false, "4+20=", 2) | Out-Null
$d.Content.Find.Execute("17+67=", $true, $false, $false, $false, $false, $true, 1, $false, "86-59=", 2) | Out-Null
$d.Content.Find.Execute("9+10=", $true, $false, $false, $false, $false, $true, 1, $false, "98-78=", 2) | Out-Null
$d.Content.Find.Execute("1+67=", $true, $false, $false, $false, $false, $true, 1, $false, "86-64=", 2) | Out-Null
$d.Content.Find.Execute("73+7=", $true, $false, $false, $false, $false, $true, 1, $false, "39+5=", 2) | Out-Null
$d.Content.Find.Execute("71-41=", $true, $false, $false, $false, $false, $true, 1, $false, "55-43=", 2) | Out-Null
$d.Content.Find.Execute("68-30=", $true, $false, $false, $false, $false, $true, 1, $false, "67+13=", 2) | Out-Null
$d.Content.Find.Execute("8+25=", $true, $false, $false, $false, $false, $true, 1, $false, "74-34=", 2) | Out-Null
$d.Content.Find.Execute("38+37=", $true, $false, $false, $false, $false, $true, 1, $false, "1+60=", 2) | Out-Null
$d.Content.Find.Execute("60+12=", $true, $false, $false, $false, $false, $true, 1, $false, "77-9=", 2) | Out-Null
$d.Content.Find.Execute("21+37=", $true, $false, $false, $false, $false, $true, 1, $false, "30-5=", 2) | Out-Null
$d.Content.Find.Execute("50+45=", $true, $false, $false, $false, $false, $true, 1, $false, "50+42=", 2) | Out-Null
$d.Content.Find.Execute("10+58=", $true, $false, $false, $false, $false, $true, 1, $false, "82-14=", 2) | Out-Null
$d.Content.Find.Execute("8+91=", $true, $false, $false, $false, $false, $true, 1, $false, "66-20=", 2) | Out-Null
$d.Content.Find.Execute("56+5=", $true, $false, $false, $false, $false, $true, 1, $false, "59+19=", 2) | Out-Null
$d.Content.Find.Execute("15+54=", $true, $false, $false, $false, $false, $true, 1, $false, "69-51=", 2) | Out-Null
$d.Content.Find.Execute("43-32=", $true, $false, $false, $false, $false, $true, 1, $false, "96-81=", 2) | Out-Null
$d.Content.Find.Execute("41+35=", $true, $false, $false, $false, $false, $true, 1, $false, "24+0=", 2) | Out-Null
$d.Content.Find.Execute("84+10=", $true, $false, $false, $false, $false, $true, 1, $false, "90-80=", 2) | Out-Null
$d.Content.Find.Execute("3+93=", $true, $false, $false, $false, $false, $true, 1, $false, "72-43=", 2) | Out-Null
$d.Content.Find.Execute("67+5=", $true, $false, $false, $false, $false, $true, 1, $false, "89-46=", 2) | Out-Null
$d.Content.Find.Execute("57-49=", $true, $false, $false, $false, $false, $true, 1, $false, "22+22=", 2) | Out-Null
$d.Content.Find.Execute("2+41=", $true, $false, $false, $false, $false, $true, 1, $false, "31+39=", 2) | Out-Null
$d.Content.Find.Execute("45-40=", $true, $false, $false, $false, $false, $true, 1, $false, "45-41=", 2) | Out-Null
$d.Content.Find.Execute("93-60=", $true, $false, $false, $false, $false, $true, 1, $false, "18+33=", 2) | Out-Null
$d.Content.Find.Execute("40+22=", $true, $false, $false, $false, $false, $true, 1, $false, "43+22=", 2) | Out-Null
$d.Content.Find.Execute("75+19=", $true, $false, $false, $false, $false, $true, 1, $false, "45-7=", 2) | Out-Null
$d.Content.Find.Execute("37+59=", $true, $false, $false, $false, $false, $true, 1, $false, "83-81=", 2) | Out-Null
$d.Content.Find.Execute("43+55=", $true, $false, $false, $false, $false, $true, 1, $false, "30+18=", 2) | Out-Null
$d.Content.Find.Execute("99-60=", $true, $false, $false, $false, $false, $true, 1, $false, "10+50=", 2) | Out-Null
$d.Content.Find.Execute("45-3=", $true, $false, $false, $false, $false, $true, 1, $false, "60-12=", 2) | Out-Null
$d.Content.Find.Execute("38-8=", $true, $false, $false, $false, $false, $true, 1, $false, "64-61=", 2) | Out-Null
$d.Content.Find.Execute("2+76=", $true, $false, $false, $false, $false, $true, 1, $false, "82-31=", 2) | Out-Null
$d.Content.Find.Execute("67-42=", $true, $false, $false, $false, $false, $true, 1, $false, "47+31=", 2) | Out-Null
$d.Content.Find.Execute("14+81=", $true, $false, $false, $false, $false, $true, 1, $false, "48+46=", 2) | Out-Null
$d.Content.Find.Execute("5+89=", $true, $false, $false, $false, $false, $true, 1, $false, "67-55=", 2) | Out-Null
$d.Content.Find.Execute("55+14=", $true, $false, $false, $false, $false, $true, 1, $false, "15+16=", 2) | Out-Null
$d.Content.Find.Execute("17+79=", $true, $false, $false, $false, $false, $true, 1, $false, "59+6=", 2) | Out-Null
$d.Content.Find.Execute("88-15=", $true, $false, $false, $false, $false, $true, 1, $false, "17-12=", 2) | Out-Null
$d.Content.Find.Execute("16+70=", $true, $false, $false, $false, $false, $true, 1, $false, "9+28=", 2) | Out-Null
$d.Content.Find.Execute("77-44=", $true, $false, $false, $false, $false, $true, 1, $false, "84-29=", 2) | Out-Null
$d.Content.Find.Execute("49-17=", $true, $false, $false, $false, $false, $true, 1, $false, "19-7=", 2) | Out-Null
$d.Content.Find.Execute("94-93=", $true, $false, $false, $false, $false, $true, 1, $false, "62+16=", 2) | Out-Null
$d.Content.Find.Execute("44-41=", $true, $false, $false, $false, $false, $true, 1, $false, "25+51=", 2) | Out-Null
$d.Content.Find.Execute("88-16=", $true, $false, $false, $false, $false, $true, 1, $false, "51-32=", 2) | Out-Null
$d.Content.Find.Execute("95-45=", $true, $false, $false, $false, $false, $true, 1, $false, "44-1=", 2) | Out-Null
$d.Content.Find.Execute("90-53=", $true, $false, $false, $false, $false, $true, 1, $false, "22+72=", 2) | Out-Null
$d.Content.Find.Execute("1+44=", $true, $false, $false, $false, $false, $true, 1, $false, "22+27=", 2) | Out-Null
$d.Content.Find.Execute("63-62=", $true, $false, $false, $false, $false, $true, 1, $false, "71+2=", 2) | Out-Null
$d.Content.Find.Execute("62-44=", $true, $false, $false, $false, $false, $true, 1, $false, "81+2=", 2) | Out-Null
$d.Content.Find.Execute("44-35=", $true, $false, $false, $false, $false, $true, 1, $false, "83-78=", 2) | Out-Null
$d.Content.Find.Execute("24-5=", $true, $false, $false, $false, $false, $true, 1, $false, "83-33=", 2) | Out-Null
$d.Content.Find.Execute("89-69=", $true, $false, $false, $false, $false, $true, 1, $false, "94-30=", 2) | Out-Null
$d.Content.Find.Execute("70-22=", $true, $false, $false, $false, $false, $true, 1, $false, "58+39=", 2) | Out-Null
$d.Content.Find.Execute("93-20=", $true, $false, $false, $false, $false, $true, 1, $false, "19+40=", 2) | Out-Null
$d.Content.Find.Execute("94-53=", $true, $false, $false, $false, $false, $true, 1, $false, "34+55=", 2) | Out-Null
$d.Content.Find.Execute("74-66=", $true, $false, $false, $false, $false, $true, 1, $false, "51-5=", 2) | Out-Null
$d.Content.Find.Execute("13+63=", $true, $false, $false, $false, $false, $true, 1, $false, "15+7=", 2) | Out-Null
$d.Content.Find.Execute("6+62=", $true, $false, $false, $false, $false, $true, 1, $false, "57-26=", 2) | Out-Null
$d.Content.Find.Execute("99-56=", $true, $false, $false, $false, $false, $true, 1, $false, "72-40=", 2) | Out-Null
$d.Content.Find.Execute("94-2=", $true, $false, $false, $false, $false, $true, 1, $false, "0+67=", 2) | Out-Null
$d.Content.Find.Execute("78+9=", $true, $false, $false, $false, $false, $true, 1, $false, "87-39=", 2) | Out-Null
$d.Content.Find.Execute("20+67=", $true, $false, $false, $false, $false, $true, 1, $false, "75-30=", 2) | Out-Null
$d.Content.Find.Execute("43+48=", $true, $false, $false, $false, $false, $true, 1, $false, "20+23=", 2) | Out-Null
$d.Content.Find.Execute("54-20=", $true, $false, $false, $false, $false, $true, 1, $false, "65+11=", 2) | Out-Null
$d.Content.Find.Execute("1+45=", $true, $false, $false, $false, $false, $true, 1, $false, "3+15=", 2) | Out-Null
$d.Content.Find.Execute("17+13=", $true, $false, $false, $false, $false, $true, 1, $false, "85-80=", 2) | Out-Null
$d.Content.Find.Execute("26+13=", $true, $false, $false, $false, $false, $true, 1, $false, "83-82=", 2) | Out-Null
$d.Content.Find.Execute("65-53=", $true, $false, $false, $false, $false, $true, 1, $false, "78-74=", 2) | Out-Null
$d.Content.Find.Execute("42+33=", $true, $false, $false, $false, $false, $true, 1, $false, "91-32=", 2) | Out-Null
